$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matches source inline-string cells)
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D12", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values from the refreshed cryptos feed
$ws.Range("D2").Value = "30.388.38"
$ws.Range("E2").Value = "  +2.63%  "
$ws.Range("D3").Value = "2.102.87"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.72%  "
$ws.Range("D5").Value = "343.89"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("D7").Value = "0.5296"
$ws.Range("E7").Value = "  +2.43%  "
$ws.Range("D8").Value = "0.4424"
$ws.Range("E8").Value = "  +1.46%  "
$ws.Range("D9").Value = "54.96"
$ws.Range("E9").Value = "  +4.11%  "
$ws.Range("D10").Value = "0.09390"
$ws.Range("E10").Value = "  +1.99%  "
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("D12").Value = "24.78"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "2.172.97"
$ws.Range("E13").Value = "  +3.40%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "8.554"
$ws.Range("E14").Value = "  +4.70%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "6.902"
$ws.Range("E15").Value = "  +2.61%  "
$ws.Range("D16").Value = "101.68"
$ws.Range("E16").Value = "  +1.97%  "
$ws.Range("D17").Value = "0.00001159"
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("D18").Value = "1.003"
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("D19").Value = "21.18"
$ws.Range("E19").Value = "  +2.10%  "
$ws.Range("D20").Value = "0.06711"
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("D21").Value = "6.338"
$ws.Range("E21").Value = "  +2.66%  "
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").Value = "  -0.59%  "
$ws.Range("D23").Value = "30.431.79"
$ws.Range("E23").Value = "  +2.63%  "
$ws.Range("D24").Value = "12.51"
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("D25").Value = "2.319"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "6.979"
$ws.Range("E26").Value = "  +11.49%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "21.88"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "162.51"
$ws.Range("E28").Value = "  +1.21%  "
$ws.Range("D29").Value = "2.508"
$ws.Range("E29").Value = "  +0.63%  "
$ws.Range("D30").Value = "134.11"
$ws.Range("E30").Value = "  +0.79%  "
$ws.Range("D31").Value = "1.137"
$ws.Range("D32").Value = "1.677"
$ws.Range("E32").Value = "  +3.40%  "
$ws.Range("D33").Value = "0.1057"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("D34").Value = "6.253"
$ws.Range("E34").Value = "  +1.55%  "
$ws.Range("D35").Value = "3.865"
$ws.Range("E35").Value = "  -2.46%  "
$ws.Range("D36").Value = "10.13"
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("D37").Value = "0.02645"
$ws.Range("E37").Value = "  +2.99%  "
$ws.Range("D38").Value = "0.06801"
$ws.Range("E38").Value = "  +1.73%  "
$ws.Range("D39").Value = "12.72"
$ws.Range("E39").Value = "  +2.11%  "
$ws.Range("D40").Value = "0.7019"
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("D41").Value = "1.344"
$ws.Range("E41").Value = "  +2.47%  "
$ws.Range("D42").Value = "0.2217"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "0.6818"
$ws.Range("E43").Value = "  -0.96%  "
$ws.Range("D44").Value = "14.51"
$ws.Range("E44").Value = "  +2.23%  "
$ws.Range("D45").Value = "2.328"
$ws.Range("E46").Value = "  -0.81%  "
$ws.Range("D47").Value = "1.353"
$ws.Range("E47").Value = "  +15.60%  "
$ws.Range("D48").Value = "3.645"
$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("D49").Value = "0.00000000357"
$ws.Range("E49").Value = "  +3.18%  "
$ws.Range("D50").Value = "1.215"
$ws.Range("E50").Value = "  +7.89%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.07349"
$ws.Range("E51").Value = "  +4.72%  "
